# Reverting Y axis: put (0, 0) at bottom-left instead of top-left.
# The sheet encodes a little pixel-art "curve" using shared strings
# "X" (mark) and "V" placed on a 15x15 (A1:O15) grid. This updates the
# marked cells to the new orientation: clearing the cells that should
# no longer be marked, and marking the cells that should now be marked.
#
# NOTE: this runtime's Range object does not reliably apply .Value to
# every area of a multi-area ("A1,B2,C3") Union range (only the first
# and last areas get written), so each cell is addressed individually
# to guarantee a correct, unambiguous result.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells that lose their mark (go back to blank), keeping their border style.
$ws.Range("B2").ClearContents()
$ws.Range("O8").ClearContents()
$ws.Range("N9").ClearContents()
$ws.Range("O9").ClearContents()

# Cells that become marked with "X" (reuses the existing shared string).
$ws.Range("M2").Value = "X"
$ws.Range("C3").Value = "X"
$ws.Range("N4").Value = "X"
$ws.Range("J5").Value = "X"
$ws.Range("K5").Value = "X"
$ws.Range("N5").Value = "X"
$ws.Range("B6").Value = "X"
$ws.Range("C6").Value = "X"
$ws.Range("J6").Value = "X"
$ws.Range("K6").Value = "X"
$ws.Range("N6").Value = "X"
$ws.Range("B7").Value = "X"
$ws.Range("K7").Value = "X"
$ws.Range("B8").Value = "X"
$ws.Range("C8").Value = "X"
$ws.Range("K8").Value = "X"
$ws.Range("E9").Value = "X"
$ws.Range("F9").Value = "X"
$ws.Range("G9").Value = "X"
$ws.Range("J9").Value = "X"
$ws.Range("K9").Value = "X"
$ws.Range("E10").Value = "X"
$ws.Range("J10").Value = "X"
$ws.Range("K10").Value = "X"
$ws.Range("L10").Value = "X"
$ws.Range("M10").Value = "X"
$ws.Range("A11").Value = "X"
$ws.Range("D11").Value = "X"
$ws.Range("E11").Value = "X"
$ws.Range("F11").Value = "X"
$ws.Range("G11").Value = "X"
$ws.Range("J11").Value = "X"
$ws.Range("K11").Value = "X"
$ws.Range("A12").Value = "X"
$ws.Range("D12").Value = "X"
$ws.Range("E12").Value = "X"
$ws.Range("N12").Value = "X"
$ws.Range("A13").Value = "X"
$ws.Range("B13").Value = "X"
$ws.Range("M13").Value = "X"
$ws.Range("N13").Value = "X"
$ws.Range("A14").Value = "X"
$ws.Range("B14").Value = "X"
$ws.Range("H14").Value = "X"
$ws.Range("I14").Value = "X"
$ws.Range("J14").Value = "X"
$ws.Range("K14").Value = "X"
$ws.Range("L14").Value = "X"
$ws.Range("M14").Value = "X"
$ws.Range("A15").Value = "X"
$ws.Range("B15").Value = "X"
$ws.Range("M15").Value = "X"

# Cell that becomes marked with "V" (the other existing shared string).
$ws.Range("N15").Value = "V"

# Reset the active selection back to the top-left cell.
$ws.Range("A1").Select() | Out-Null
